# S-01015: se cargan las horas insumidas para parsear el archivo HF.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Horas insumidas")

# Insert a new row above the existing row 44 (2010-10-01, Duilio / Scrolling en
# ABM de campaña / S-01018 / 1h), shifting the rest of the table (including the
# "Total Sprint 2" footer row) down by one.
$ws.Rows("44:44").Insert()

# The inserted row picks up the border/fill from the "Total Sprint 1" banner
# row above it; reset it back to the plain look used by the other data rows.
$ws.Range("B44:F44").ClearFormats()

$ws.Range("B44").Value = "9/29/2010"
$ws.Range("B44").NumberFormat = "d-mmm"
$ws.Range("C44").Value = "Duilio"
$ws.Range("E44").Value = "S-01015"
$ws.Range("D44").Value = "Proceso archivo HF"
$ws.Range("F44").Value = 4

# Keep the view roughly where the author left it.
$ws.Range("E46").Select()
